# Auto-generated edit script applying the Typhon_Profits.xlsx value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1000.9091
$ws.Range("I58").Value = 117.36842
$ws.Range("J58").Value = 6596.6665
$ws.Range("K58").Value = 352.10526
$ws.Range("L58").Value = 19789.9995
$ws.Range("M58").Value = -202.10526
$ws.Range("N58").Value = -20089.9995
$ws.Range("H86").Value = 27884
$ws.Range("I86").Value = 3525
$ws.Range("J86").Value = 52243
$ws.Range("K86").Value = 3525
$ws.Range("L86").Value = 52243
$ws.Range("M86").Value = -2402
$ws.Range("N86").Value = -54489
$ws.Range("H87").Value = 40354
$ws.Range("J87").Value = 40354
$ws.Range("L87").Value = 40354
$ws.Range("N87").Value = -42850
$ws.Range("H88").Value = 2172
$ws.Range("J88").Value = 2234
$ws.Range("L88").Value = 2234
$ws.Range("N88").Value = -3046
$ws.Range("H89").Value = 27884
$ws.Range("I89").Value = 3525
$ws.Range("J89").Value = 52243
$ws.Range("K89").Value = 17625
$ws.Range("L89").Value = 261215
$ws.Range("M89").Value = -12009
$ws.Range("N89").Value = -272447
$ws.Range("H90").Value = 40354
$ws.Range("J90").Value = 40354
$ws.Range("L90").Value = 121062
$ws.Range("N90").Value = -133542
$ws.Range("H91").Value = 2172
$ws.Range("J91").Value = 2234
$ws.Range("L91").Value = 2234
$ws.Range("N91").Value = -5042
$ws.Range("H107").Value = 1032.7693
$ws.Range("I107").Value = 1087.375
$ws.Range("J107").Value = 945.4
$ws.Range("K107").Value = 1087.375
$ws.Range("L107").Value = 945.4
$ws.Range("M107").Value = 832.625
$ws.Range("N107").Value = -4785.4
$ws.Range("H129").Value = 952.80853
$ws.Range("J129").Value = 1073.5526
$ws.Range("L129").Value = 3220.6578
$ws.Range("N129").Value = -13220.6578
$ws.Range("H132").Value = 26680.342
$ws.Range("I132").Value = 31891.117
$ws.Range("J132").Value = 1370.8572
$ws.Range("K132").Value = 95673.351
$ws.Range("L132").Value = 4112.571599999999
$ws.Range("M132").Value = -93143.351
$ws.Range("N132").Value = -9172.571599999999
$ws.Range("H137").Value = 18153.361
$ws.Range("I137").Value = 1782.825
$ws.Range("J137").Value = 49335.332
$ws.Range("K137").Value = 5348.475
$ws.Range("L137").Value = 148005.996
$ws.Range("M137").Value = -2798.475
$ws.Range("N137").Value = -153105.996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16599.648
$ws.Range("I32").Value = 17831.629
$ws.Range("K32").Value = 17831.629
$ws.Range("M32").Value = -17544.629
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H102").Value = 4478.8887
$ws.Range("I102").Value = 2901.4285
$ws.Range("K102").Value = 2901.4285
$ws.Range("M102").Value = -1279.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12795.419
$ws.Range("I31").Value = 38987.625
$ws.Range("J31").Value = 3685.087
$ws.Range("K31").Value = 38987.625
$ws.Range("L31").Value = 3685.087
$ws.Range("M31").Value = -38692.625
$ws.Range("N31").Value = -4275.087
$ws.Range("H34").Value = 12795.419
$ws.Range("I34").Value = 38987.625
$ws.Range("J34").Value = 3685.087
$ws.Range("K34").Value = 38987.625
$ws.Range("L34").Value = 3685.087
$ws.Range("M34").Value = -38785.625
$ws.Range("N34").Value = -4089.087
$ws.Range("H132").Value = 18677
$ws.Range("I132").Value = 24816.523
$ws.Range("K132").Value = 74449.569
$ws.Range("M132").Value = -71919.569

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 7165.143
$ws.Range("I2").Value = 14312.714
$ws.Range("J2").Value = 17.571428
$ws.Range("K2").Value = 85876.284
$ws.Range("L2").Value = 105.428568
$ws.Range("M2").Value = -85763.284
$ws.Range("N2").Value = -331.428568
$ws.Range("H5").Value = 736.5806
$ws.Range("I5").Value = 622.9
$ws.Range("J5").Value = 790.7143
$ws.Range("K5").Value = 1868.7
$ws.Range("L5").Value = 2372.1429
$ws.Range("M5").Value = -1756.7
$ws.Range("N5").Value = -2596.1429
$ws.Range("H17").Value = 457.23077
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 616
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 1848
$ws.Range("M17").Value = -131
$ws.Range("N17").Value = -2186
$ws.Range("H34").Value = 742.875
$ws.Range("J34").Value = 829
$ws.Range("L34").Value = 2487
$ws.Range("N34").Value = -2655
$ws.Range("H39").Value = 1976
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 1976
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 5928
$ws.Range("N39").Value = -6516
$ws.Range("M39").ClearContents()
$ws.Range("H55").Value = 2300
$ws.Range("J55").Value = 2300
$ws.Range("L55").Value = 6900
$ws.Range("N55").Value = -7254
$ws.Range("H131").Value = 115755.2
$ws.Range("J131").Value = 124273.87
$ws.Range("L131").Value = 372821.61
$ws.Range("N131").Value = -382901.61
$ws.Range("H135").Value = 736.5806
$ws.Range("I135").Value = 622.9
$ws.Range("J135").Value = 790.7143
$ws.Range("K135").Value = 5606.099999999999
$ws.Range("L135").Value = 7116.428699999999
$ws.Range("M135").Value = -3071.099999999999
$ws.Range("N135").Value = -12186.4287
$ws.Range("H137").Value = 9420.823
$ws.Range("J137").Value = 3759.6875
$ws.Range("L137").Value = 11279.0625
$ws.Range("N137").Value = -21479.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8650.526
$ws.Range("I80").Value = 13728.889
$ws.Range("J80").Value = 4080
$ws.Range("K80").Value = 13728.889
$ws.Range("L80").Value = 4080
$ws.Range("M80").Value = -12730.889
$ws.Range("N80").Value = -6076
$ws.Range("H83").Value = 8650.526
$ws.Range("I83").Value = 13728.889
$ws.Range("J83").Value = 4080
$ws.Range("K83").Value = 68644.44499999999
$ws.Range("L83").Value = 20400
$ws.Range("M83").Value = -63652.44499999999
$ws.Range("N83").Value = -30384
$ws.Range("H102").Value = 1481.3478
$ws.Range("I102").Value = 1481.8889
$ws.Range("J102").Value = 1479.4
$ws.Range("K102").Value = 1481.8889
$ws.Range("L102").Value = 1479.4
$ws.Range("M102").Value = 140.1111000000001
$ws.Range("N102").Value = -4723.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4262.2144
$ws.Range("I7").Value = 4410
$ws.Range("K7").Value = 4410
$ws.Range("M7").Value = -4298
$ws.Range("H22").Value = 1696.3334
$ws.Range("I22").Value = 3592
$ws.Range("J22").Value = 748.5
$ws.Range("K22").Value = 3592
$ws.Range("L22").Value = 748.5
$ws.Range("M22").Value = -3297
$ws.Range("N22").Value = -1338.5
$ws.Range("H27").Value = 1696.3334
$ws.Range("I27").Value = 3592
$ws.Range("J27").Value = 748.5
$ws.Range("K27").Value = 3592
$ws.Range("L27").Value = 748.5
$ws.Range("M27").Value = -3485
$ws.Range("N27").Value = -962.5
$ws.Range("H32").Value = 3275.3333
$ws.Range("I32").Value = 3275.3333
$ws.Range("K32").Value = 3275.3333
$ws.Range("M32").Value = -2958.3333
$ws.Range("H126").Value = 4262.2144
$ws.Range("I126").Value = 4410
$ws.Range("K126").Value = 13230
$ws.Range("M126").Value = -10760

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3041.4
$ws.Range("J15").Value = 3041.4
$ws.Range("L15").Value = 3041.4
$ws.Range("N15").Value = -3617.4
$ws.Range("H54").Value = 14833.333
$ws.Range("J54").Value = 14833.333
$ws.Range("L54").Value = 14833.333
$ws.Range("N54").Value = -15873.333
$ws.Range("H81").Value = 1733.091
$ws.Range("I81").Value = 1311.2222
$ws.Range("J81").Value = 3631.5
$ws.Range("K81").Value = 2622.4444
$ws.Range("L81").Value = 7263
$ws.Range("M81").Value = -1561.4444
$ws.Range("N81").Value = -9385
$ws.Range("H84").Value = 1733.091
$ws.Range("I84").Value = 1311.2222
$ws.Range("J84").Value = 3631.5
$ws.Range("K84").Value = 13112.222
$ws.Range("L84").Value = 36315
$ws.Range("M84").Value = -7808.222
$ws.Range("N84").Value = -46923
$ws.Range("H122").Value = 1362.5161
$ws.Range("I122").Value = 1234.3846
$ws.Range("K122").Value = 3703.1538
$ws.Range("M122").Value = -1253.1538
$ws.Range("H132").Value = 1934.1562
$ws.Range("I132").Value = 1753.9584
$ws.Range("J132").Value = 2474.75
$ws.Range("K132").Value = 5261.8752
$ws.Range("L132").Value = 7424.25
$ws.Range("M132").Value = -2731.8752
$ws.Range("N132").Value = -12484.25

